$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "25.908.85"
Set-TextValue "E2" "  +0.33%  "

Set-TextValue "D3" "1.642.17"
Set-TextValue "E3" "  +0.46%  "

Set-TextValue "D4" "1.007"
Set-TextValue "E4" "  +0.47%  "

Set-TextValue "D5" "215.73"
Set-TextValue "E5" "  +0.30%  "

Set-TextValue "D6" "0.5064"
Set-TextValue "E6" "  +0.92%  "

Set-TextValue "E7" "  +0.38%  "

Set-TextValue "D8" "0.2579"
Set-TextValue "E8" "  +0.38%  "

Set-TextValue "D9" "0.06429"
Set-TextValue "E9" "  +0.34%  "

Set-TextValue "D10" "19.75"
Set-TextValue "E10" "  +0.78%  "

Set-TextValue "D11" "0.07798"
Set-TextValue "E11" "  +1.71%  "

Set-TextValue "D12" "4.309"
Set-TextValue "E12" "  +1.78%  "

Set-TextValue "D13" "1.653.96"
Set-TextValue "E13" "  +1.12%  "

Set-TextValue "D14" "0.5465"
Set-TextValue "E14" "  +0.20%  "

Set-TextValue "D15" "0.0₅7908"
Set-TextValue "E15" "  -0.08%  "

Set-TextValue "D16" "65.24"
Set-TextValue "E16" "  +2.76%  "

Set-TextValue "D17" "26.002.31"
Set-TextValue "E17" "  +0.62%  "

Set-TextValue "E18" "  +0.44%  "

Set-TextValue "D19" "198.51"
Set-TextValue "E19" "  -2.22%  "

Set-TextValue "D20" "4.417"
Set-TextValue "E20" "  +2.67%  "

Set-TextValue "D21" "10.01"
Set-TextValue "E21" "  +0.82%  "

Set-TextValue "D22" "6.047"
Set-TextValue "E22" "  +1.36%  "

Set-TextValue "D23" "1.009"
Set-TextValue "E23" "  +0.56%  "

Set-TextValue "D24" "1.873"
Set-TextValue "E24" "  -3.21%  "

Set-TextValue "D25" "140.95"
Set-TextValue "E25" "  -0.02%  "

Set-TextValue "D26" "0.1149"
Set-TextValue "E26" "  +0.46%  "

Set-TextValue "D27" "6.903"
Set-TextValue "E27" "  +3.21%  "

Set-TextValue "D28" "15.74"
Set-TextValue "E28" "  +0.19%  "

Set-TextValue "D29" "1.244"
Set-TextValue "E29" "  +0.55%  "

Set-TextValue "D30" "0.05052"
Set-TextValue "E30" "  +1.51%  "

Set-TextValue "D31" "3.277"
Set-TextValue "E31" "  +0.07%  "

Set-TextValue "D32" "3.205"
Set-TextValue "E32" "  +0.69%  "

Set-TextValue "D33" "1.542"
Set-TextValue "E33" "  +0.70%  "

Set-TextValue "D34" "2.375"
Set-TextValue "E34" "  +1.00%  "

Set-TextValue "D35" "0.8961"
Set-TextValue "E35" "  +0.52%  "

Set-TextValue "D36" "2.600"
Set-TextValue "E36" "  -0.88%  "

Set-TextValue "D37" "1.134.98"
Set-TextValue "E37" "  -3.33%  "

Set-TextValue "D38" "0.5546"
Set-TextValue "E38" "  -0.50%  "

Set-TextValue "D39" "0.01566"
Set-TextValue "E39" "  +0.55%  "

Set-TextValue "D40" "1.010"
Set-TextValue "E40" "  +0.72%  "

Set-TextValue "D41" "5.700"
Set-TextValue "E41" "  +1.14%  "

Set-TextValue "D42" "0.8179"
Set-TextValue "E42" "  +1.96%  "

Set-TextValue "D43" "99.86"
Set-TextValue "E43" "  +0.56%  "

Set-TextValue "E44" "  +6.76%  "

Set-TextValue "D45" "1.778.35"
Set-TextValue "E45" "  +0.40%  "

Set-TextValue "D46" "0.4539"
Set-TextValue "E46" "  +0.67%  "

Set-TextValue "D47" "55.41"
Set-TextValue "E47" "  +1.13%  "

Set-TextValue "D48" "1.006"
Set-TextValue "E48" "  +0.16%  "

Set-TextValue "D49" "0.05089"
Set-TextValue "E49" "  +1.15%  "

Set-TextValue "B50" "Algorand"
Set-TextValue "C50" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D50" "0.09596"
Set-TextValue "E50" "  +3.63%  "

Set-TextValue "B51" "USDD"
Set-TextValue "C51" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue "D51" "1.007"
Set-TextValue "E51" "  +0.50%  "
